$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Needleman-Wunsch traceback matrix (rows 8-11, cols E:G) ---
# Replace the plain numeric scores with direction/score labels ("L, 1", "U, 2", "D, 0", ...)
$ws.Range("F8").Value = "L, 1"
$ws.Range("G8").Value = "L,2"

$ws.Range("E9").Value = "U, 1"
$ws.Range("F9").Value = "D, 1"
$ws.Range("G9").Value = "D, 2"

$ws.Range("E10").Value = "U, 2"
$ws.Range("F10").Value = "D, 1"
$ws.Range("G10").Value = "D, 2"

$ws.Range("E11").Value = "U, 3"
$ws.Range("F11").Value = "U, 2"
$ws.Range("G11").Value = "D, 1"

$ws.Range("E8").Value = "D, 0"

# --- Highlight the traceback path in green (new fill FF92D050) ---
$tracebackCells = "E8", "E9", "F10", "G11"
foreach ($addr in $tracebackCells) {
    $ws.Range($addr).Interior.Color = 5296274
    $ws.Range($addr).HorizontalAlignment = -4108
}

# --- View changes: hide column H, zoom to 160%, move the selection to F8 ---
$ws.Columns("H").Hidden = $true
$excel.ActiveWindow.Zoom = 160
$ws.Range("F8").Select()
